$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9005614640222286
$ws.Range("J2").Value = 0.9005614640222285
$ws.Range("M2").Value = 6.101885666666667
$ws.Range("N2").Value = 18.305657
$ws.Range("O2").Value = 0.1093737608697887
$ws.Range("P2").Value = 0.1093737608697887
$ws.Range("Q2").Value = 0.4779383306892223
$ws.Range("R2").Value = 4.301444976203
$ws.Range("S2").Value = 0.09849779421451409
$ws.Range("T2").Value = 0.09849779421451407

$ws.Range("I3").Value = 0.9005614640222286
$ws.Range("J3").Value = 0.9005614640222285
$ws.Range("N3").Value = 87.53628900000001
$ws.Range("O3").Value = 0.5230171820937495
$ws.Range("P3").Value = 0.5230171820937495
$ws.Range("Q3").Value = 2.285465516992334
$ws.Range("S3").Value = 0.4710091192151276
$ws.Range("T3").Value = 0.4710091192151275

$ws.Range("I4").Value = 0.9005614640222286
$ws.Range("J4").Value = 0.9005614640222285
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.146644
$ws.Range("N4").Value = 0.439932
$ws.Range("O4").Value = 0.002628532664354407
$ws.Range("P4").Value = 0.002628532664354407
$ws.Range("Q4").Value = 0.01148608682533333
$ws.Range("R4").Value = 0.103374781428
$ws.Range("S4").Value = 0.002367155224441254
$ws.Range("T4").Value = 0.002367155224441254

$ws.Range("I5").Value = 0.9005614640222286
$ws.Range("J5").Value = 0.9005614640222285
$ws.Range("M5").Value = 15.02284966666667
$ws.Range("N5").Value = 45.068549
$ws.Range("O5").Value = 0.2692783275177917
$ws.Range("P5").Value = 0.2692783275177917
$ws.Range("Q5").Value = 1.176684730607889
$ws.Range("R5").Value = 10.590162575471
$ws.Range("S5").Value = 0.2425016848588796
$ws.Range("T5").Value = 0.2425016848588797

$ws.Range("I6").Value = 0.9005614640222286
$ws.Range("J6").Value = 0.9005614640222285
$ws.Range("M6").Value = 5.288900666666667
$ws.Range("N6").Value = 15.866702
$ws.Range("O6").Value = 0.09480134312252211
$ws.Range("P6").Value = 0.09480134312252211
$ws.Range("Q6").Value = 0.4142601965842222
$ws.Range("R6").Value = 3.728341769258
$ws.Range("S6").Value = 0.08537443635369214
$ws.Range("T6").Value = 0.08537443635369213

$ws.Range("I7").Value = 0.9005614640222286
$ws.Range("J7").Value = 0.9005614640222285
$ws.Range("M7").Value = 0.050258
$ws.Range("N7").Value = 0.150774
$ws.Range("O7").Value = 0.0009008537317934847
$ws.Range("P7").Value = 0.0009008537317934848
$ws.Range("Q7").Value = 0.003936524860666667
$ws.Range("R7").Value = 0.035428723746
$ws.Range("S7").Value = 0.0008112741555738286
$ws.Range("T7").Value = 0.0008112741555738286

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.008648666666666667
$ws.Range("H8").Value = 0.025946
$ws.Range("I8").Value = 0.0994385359777714
$ws.Range("J8").Value = 0.09943853597777139
$ws.Range("M8").Value = 6.101885666666667
$ws.Range("N8").Value = 18.305657
$ws.Range("O8").Value = 0.1093737608697887
$ws.Range("P8").Value = 0.1093737608697887
$ws.Range("Q8").Value = 0.05277317516911112
$ws.Range("R8").Value = 0.474958576522
$ws.Range("S8").Value = 0.01087596665527465
$ws.Range("T8").Value = 0.01087596665527465

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.008648666666666667
$ws.Range("H9").Value = 0.025946
$ws.Range("I9").Value = 0.0994385359777714
$ws.Range("J9").Value = 0.09943853597777139
$ws.Range("N9").Value = 87.53628900000001
$ws.Range("O9").Value = 0.5230171820937495
$ws.Range("P9").Value = 0.5230171820937495
$ws.Range("Q9").Value = 0.2523573949326667
$ws.Range("R9").Value = 2.271216554394
$ws.Range("S9").Value = 0.05200806287862193
$ws.Range("T9").Value = 0.05200806287862192

$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.008648666666666667
$ws.Range("H10").Value = 0.025946
$ws.Range("I10").Value = 0.0994385359777714
$ws.Range("J10").Value = 0.09943853597777139
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.146644
$ws.Range("N10").Value = 0.439932
$ws.Range("O10").Value = 0.002628532664354407
$ws.Range("P10").Value = 0.002628532664354407
$ws.Range("Q10").Value = 0.001268275074666667
$ws.Range("R10").Value = 0.011414475672
$ws.Range("S10").Value = 0.000261377439913153
$ws.Range("T10").Value = 0.000261377439913153

$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.008648666666666667
$ws.Range("H11").Value = 0.025946
$ws.Range("I11").Value = 0.0994385359777714
$ws.Range("J11").Value = 0.09943853597777139
$ws.Range("M11").Value = 15.02284966666667
$ws.Range("N11").Value = 45.068549
$ws.Range("O11").Value = 0.2692783275177917
$ws.Range("P11").Value = 0.2692783275177917
$ws.Range("Q11").Value = 0.1299276191504445
$ws.Range("R11").Value = 1.169348572354
$ws.Range("S11").Value = 0.02677664265891204
$ws.Range("T11").Value = 0.02677664265891204

$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.008648666666666667
$ws.Range("H12").Value = 0.025946
$ws.Range("I12").Value = 0.0994385359777714
$ws.Range("J12").Value = 0.09943853597777139
$ws.Range("M12").Value = 5.288900666666667
$ws.Range("N12").Value = 15.866702
$ws.Range("O12").Value = 0.09480134312252211
$ws.Range("P12").Value = 0.09480134312252211
$ws.Range("Q12").Value = 0.04574193889911111
$ws.Range("R12").Value = 0.411677450092
$ws.Range("S12").Value = 0.009426906768829967
$ws.Range("T12").Value = 0.009426906768829965

$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.008648666666666667
$ws.Range("H13").Value = 0.025946
$ws.Range("I13").Value = 0.0994385359777714
$ws.Range("J13").Value = 0.09943853597777139
$ws.Range("M13").Value = 0.050258
$ws.Range("N13").Value = 0.150774
$ws.Range("O13").Value = 0.0009008537317934847
$ws.Range("P13").Value = 0.0009008537317934848
$ws.Range("Q13").Value = 0.0004346646893333333
$ws.Range("R13").Value = 0.003911982204
$ws.Range("S13").Value = 0.00008957957621965606
$ws.Range("T13").Value = 0.00008957957621965606
